$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(28, "Old Acad", "C01", 106),
    @(29, "Old Acad", "C02", 42),
    @(30, "Old Acad", "C03", 18),
    @(31, "Old Acad", "C11", 94),
    @(32, "Old Acad", "C12", 39),
    @(33, "Old Acad", "C13", 18),
    @(34, "Old Acad", "C21", 80),
    @(35, "Old Acad", "C22", 18),
    @(36, "Old Acad", "C24", 18),
    @(37, "LHC", "C101", 169),
    @(38, "LHC", "C102", 281),
    @(39, "LHC", "C201", 177),
    @(40, "LHC", "C208", 20),
    @(41, "LHC", "C209", 22),
    @(42, "LHC", "C210", 25),
    @(43, "Old Acad", "C211", 24),
    @(44, "LHC", "C212", 28),
    @(45, "LHC", "C213", 31),
    @(46, "LHC", "C214", 20),
    @(47, "LHC", "C215", 22),
    @(48, "LHC", "C216", 17),
    @(49, "RnD", "A006", 58),
    @(50, "RnD", "A007", 56),
    @(51, "RnD", "A106", 48),
    @(52, "RnD", "B105", 42)
)

$rowIndex = 2
foreach ($entry in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $entry[0]
    $ws.Cells.Item($rowIndex, 2).Value = $entry[1]
    $ws.Cells.Item($rowIndex, 3).Value = $entry[2]
    $ws.Cells.Item($rowIndex, 4).Value = $entry[3]
    $rowIndex++
}
